$p = $ppt.ActivePresentation

# --- 1. Bump the cached "last modified" date field text from 2020/11/25 to
#        2020/11/26 everywhere it appears (slide master, every slide layout,
#        the notes master and the handout master). ---
function Update-DateShapes {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "2020/11/25") {
                $tr.Text = "2020/11/26"
            }
        }
    }
}

Update-DateShapes $p.SlideMaster.Shapes
Update-DateShapes $p.NotesMaster.Shapes
Update-DateShapes $p.HandoutMaster.Shapes

for ($j = 1; $j -le $p.SlideMaster.CustomLayouts.Count; $j++) {
    Update-DateShapes $p.SlideMaster.CustomLayouts.Item($j).Shapes
}

# --- 2. Clear out the (now stale) speaker notes on slide 1: the notes body
#        placeholder held a long explanation of the diagram that has been
#        removed entirely, leaving an empty notes page. ---
$s = $p.Slides.Item(1)
$notes = $s.NotesPage
for ($i = 1; $i -le $notes.Shapes.Count; $i++) {
    $sh = $notes.Shapes.Item($i)
    if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
        if ($sh.Type -eq 14 -or $sh.PlaceholderFormat.Type -eq 2) {
            $sh.TextFrame.TextRange.Text = ""
        }
    }
}
